$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CGXII test")

$ws.Range("H13:M13").WrapText = $true

$ws.Range("A13").Value = "OD600 (after 19h)"
$ws.Range("B13").Value = 0.042
$ws.Range("C13").Value = 0.019
$ws.Range("D13").Value = 0.109
$ws.Range("E13").Value = 0.114
$ws.Range("F13").Value = 0.099
$ws.Range("G13").Value = 0.171
$ws.Range("H13").Value = 0.052
$ws.Range("I13").Value = 0.043
$ws.Range("J13").Value = 0.02
$ws.Range("K13").Value = 0.025
$ws.Range("L13").Value = 0.024
$ws.Range("M13").Value = 0.015
$ws.Range("N13").Value = "maybe waited too long before logphase"

$ws.Range("N13").Select()
